$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ty Farrington")
$ws.Cells.Item(30, 2).Value = 2
$ws.Cells.Item(30, 3).Value = "Finished singleplayer blackjack"
Write-Host ("B30=" + $ws.Cells.Item(30,2).Value)
Write-Host ("C30=" + $ws.Cells.Item(30,3).Value)
